$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" data column (H) next to the existing 2016-2019 columns,
# copying each row-cell number format from the corresponding column G cell
# (rows 8 and 20 intentionally pick up the "0.0" format used in column G/F of row 10).
$ws.Range("H4").Value = 2020
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H5").Value = 42.2
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H7").Value = 42.5
$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H8").Value = 42
$ws.Range("G10").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H10").Value = 50.9
$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H11").Value = 36.9
$ws.Range("G11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H12").Value = 34.799999999999997
$ws.Range("G12").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("G13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H14").Value = 30.7
$ws.Range("G14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H15").Value = 48.8
$ws.Range("G15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("G16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H17").Value = 61.1
$ws.Range("G17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H18").Value = 56.7
$ws.Range("G18").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H19").Value = 41.6
$ws.Range("G19").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H20").Value = 49
$ws.Range("G10").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H21").Value = 43.5
$ws.Range("G21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H22").Value = 33.9
$ws.Range("G22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H23").Value = 34.6
$ws.Range("G23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H24").Value = 23.6
$ws.Range("G24").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("H25").Value = 35.9
$ws.Range("G25").Copy()
$ws.Range("H25").PasteSpecial(-4122)

# Clear clipboard/marching-ants mode left over from the Copy() calls above.
$excel.CutCopyMode = $false

# Update the saved view: scroll position resets and the active selection moves to B13.
$ws.Range("B13").Select()
